$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the data through 02/05 (new rows 239-244), matching the formatting
# (date style) already used in column A by copying it down from the last
# existing data row before filling in the values.
$ws.Range("A238").Copy()
$ws.Range("A239:A244").PasteSpecial(-4122)

$data = @(
    @(44313, 1, 29, 287.9841112214499),
    @(44314, 0, 23, 228.4011916583913),
    @(44315, 5, 24, 238.3316782522344),
    @(44316, 9, 31, 307.8450844091361),
    @(44317, 7, 34, 337.6365441906653),
    @(44318, 5, 34, 337.6365441906653)
)

$startRow = 239
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
